$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 83.07291666666666
$ws.Range("C2").Value = 80.72916666666666
$ws.Range("D2").Value = 57.03125
$ws.Range("E2").Value = 52.60416666666667
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 47.91666666666666
$ws.Range("H2").Value = 47.39583333333334
$ws.Range("I2").Value = 47.39583333333334
$ws.Range("J2").Value = 47.13541666666666
$ws.Range("K2").Value = 47.13541666666666
$ws.Range("L2").Value = 47.13541666666666
$ws.Range("M2").Value = 47.13541666666666
$ws.Range("N2").Value = 47.13541666666666
$ws.Range("O2").Value = 47.13541666666666
$ws.Range("P2").Value = 47.13541666666666
$ws.Range("Q2").Value = 47.13541666666666
$ws.Range("R2").Value = 47.13541666666666
$ws.Range("S2").Value = 47.13541666666666
$ws.Range("T2").Value = 47.13541666666666
$ws.Range("U2").Value = 47.13541666666666
